$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.125615"
$ws.Range("H2").Value = [double]"0.25123"
$ws.Range("I2").Value = [double]"0.02647478672532295"
$ws.Range("J2").Value = [double]"0.01780700335556722"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.2272265"
$ws.Range("N2").Value = [double]"0.454453"
$ws.Range("O2").Value = [double]"0.08704083604617911"
$ws.Range("P2").Value = [double]"0.08229687998280369"
$ws.Range("Q2").Value = [double]"0.0285430567975"
$ws.Range("R2").Value = [double]"0.11417222719"
$ws.Range("S2").Value = [double]"0.002304387570716394"
$ws.Range("T2").Value = [double]"0.001465460818006498"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.125615"
$ws.Range("H3").Value = [double]"0.25123"
$ws.Range("I3").Value = [double]"0.02647478672532295"
$ws.Range("J3").Value = [double]"0.01780700335556722"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"0.2747246666666667"
$ws.Range("N3").Value = [double]"0.824174"
$ws.Range("O3").Value = [double]"0.1052353694185077"
$ws.Range("P3").Value = [double]"0.149249644656207"
$ws.Range("Q3").Value = [double]"0.03450953900333333"
$ws.Range("R3").Value = [double]"0.20705723402"
$ws.Range("S3").Value = [double]"0.002786083961315565"
$ws.Range("T3").Value = [double]"0.002657688923210293"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.125615"
$ws.Range("H4").Value = [double]"0.25123"
$ws.Range("I4").Value = [double]"0.02647478672532295"
$ws.Range("J4").Value = [double]"0.01780700335556722"
$ws.Range("K4").Value = [double]"2"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"2.082377"
$ws.Range("N4").Value = [double]"4.164754"
$ws.Range("O4").Value = [double]"0.7976703203338269"
$ws.Range("P4").Value = [double]"0.7541951755096822"
$ws.Range("Q4").Value = [double]"0.261577786855"
$ws.Range("R4").Value = [double]"1.04631114742"
$ws.Range("S4").Value = [double]"0.02111815160795811"
$ws.Range("T4").Value = [double]"0.01342995602105352"

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = [double]"2"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.125615"
$ws.Range("H5").Value = [double]"0.25123"
$ws.Range("I5").Value = [double]"0.02647478672532295"
$ws.Range("J5").Value = [double]"0.01780700335556722"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"0.02610733333333333"
$ws.Range("N5").Value = [double]"0.078322"
$ws.Range("O5").Value = [double]"0.0100006122537187"
$ws.Range("P5").Value = [double]"0.01418332860386696"
$ws.Range("Q5").Value = [double]"0.003279472676666667"
$ws.Range("R5").Value = [double]"0.01967683606"
$ws.Range("S5").Value = [double]"0.0002647640765398541"
$ws.Range("T5").Value = [double]"0.0002525625800421715"

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = [double]"2"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"0.125615"
$ws.Range("H6").Value = [double]"0.25123"
$ws.Range("I6").Value = [double]"0.02647478672532295"
$ws.Range("J6").Value = [double]"0.01780700335556722"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.3333333333333333"
$ws.Range("M6").Value = [double]"0.000138"
$ws.Range("N6").Value = [double]"0.000414"
$ws.Range("O6").Value = [double]"5.286194776741585E-05"
$ws.Range("P6").Value = [double]"7.49712474400669E-05"
$ws.Range("Q6").Value = [double]"1.733487E-05"
$ws.Range("R6").Value = [double]"0.00010400922"
$ws.Range("S6").Value = [double]"1.399508793027497E-06"
$ws.Range("T6").Value = [double]"1.335013254736332E-06"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"4.619088000000001"
$ws.Range("H7").Value = [double]"13.857264"
$ws.Range("I7").Value = [double]"0.9735252132746771"
$ws.Range("J7").Value = [double]"0.9821929966444328"
$ws.Range("K7").Value = [double]"2"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"0.2272265"
$ws.Range("N7").Value = [double]"0.454453"
$ws.Range("O7").Value = [double]"0.08704083604617911"
$ws.Range("P7").Value = [double]"0.08229687998280369"
$ws.Range("Q7").Value = [double]"1.049579199432"
$ws.Range("R7").Value = [double]"6.297475196592"
$ws.Range("S7").Value = [double]"0.08473644847546273"
$ws.Range("T7").Value = [double]"0.08083141916479719"

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"4.619088000000001"
$ws.Range("H8").Value = [double]"13.857264"
$ws.Range("I8").Value = [double]"0.9735252132746771"
$ws.Range("J8").Value = [double]"0.9821929966444328"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"0.2747246666666667"
$ws.Range("N8").Value = [double]"0.824174"
$ws.Range("O8").Value = [double]"0.1052353694185077"
$ws.Range("P8").Value = [double]"0.149249644656207"
$ws.Range("Q8").Value = [double]"1.268977411104"
$ws.Range("R8").Value = [double]"11.420796699936"
$ws.Range("S8").Value = [double]"0.1024492854571921"
$ws.Range("T8").Value = [double]"0.1465919557329967"

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"4.619088000000001"
$ws.Range("H9").Value = [double]"13.857264"
$ws.Range("I9").Value = [double]"0.9735252132746771"
$ws.Range("J9").Value = [double]"0.9821929966444328"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"2.082377"
$ws.Range("N9").Value = [double]"4.164754"
$ws.Range("O9").Value = [double]"0.7976703203338269"
$ws.Range("P9").Value = [double]"0.7541951755096822"
$ws.Range("Q9").Value = [double]"9.618682612176002"
$ws.Range("R9").Value = [double]"57.71209567305601"
$ws.Range("S9").Value = [double]"0.7765521687258689"
$ws.Range("T9").Value = [double]"0.7407652194886287"

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt5a"
$ws.Range("C10").Value = "Fzd3"
$ws.Range("D10").Value = "Neutrophils"
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"4.619088000000001"
$ws.Range("H10").Value = [double]"13.857264"
$ws.Range("I10").Value = [double]"0.9735252132746771"
$ws.Range("J10").Value = [double]"0.9821929966444328"
$ws.Range("K10").Value = [double]"2"
$ws.Range("L10").Value = [double]"0.6666666666666666"
$ws.Range("M10").Value = [double]"0.02610733333333333"
$ws.Range("N10").Value = [double]"0.078322"
$ws.Range("O10").Value = [double]"0.0100006122537187"
$ws.Range("P10").Value = [double]"0.01418332860386696"
$ws.Range("Q10").Value = [double]"0.120592070112"
$ws.Range("R10").Value = [double]"1.085328631008"
$ws.Range("S10").Value = [double]"0.009735848177178852"
$ws.Range("T10").Value = [double]"0.01393076602382479"

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt5a"
$ws.Range("C11").Value = "Fzd3"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"4.619088000000001"
$ws.Range("H11").Value = [double]"13.857264"
$ws.Range("I11").Value = [double]"0.9735252132746771"
$ws.Range("J11").Value = [double]"0.9821929966444328"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.3333333333333333"
$ws.Range("M11").Value = [double]"0.000138"
$ws.Range("N11").Value = [double]"0.000414"
$ws.Range("O11").Value = [double]"5.286194776741585E-05"
$ws.Range("P11").Value = [double]"7.49712474400669E-05"
$ws.Range("Q11").Value = [double]"0.0006374341440000001"
$ws.Range("R11").Value = [double]"0.005736907296"
$ws.Range("S11").Value = [double]"5.146243897438836E-05"
$ws.Range("T11").Value = [double]"7.363623418533058E-05"

